# Bitacora: agregar Log y Resumen de Git/Vercel + hoja Ref Git y Vercel

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Sheet "Log" - append rows 12-14
# ---------------------------------------------------------------
$log = $wb.Worksheets.Item("Log")

$logRows = @(
    @("27/02/2025", "14:00", "Repositorio Git en GitHub", "Crear repo fornitalia en GitHub (lucasbustosmartin-coder). git init, .gitignore (node_modules, .venv, .env), primer commit con dashboard, bitácora, scripts, SQL. Remote origin: https://github.com/lucasbustosmartin-coder/fornitalia.git. Push a rama main.", "Diagnostico"),
    @("27/02/2025", "14:15", "Despliegue en Vercel", "Conectar cuenta GitHub a Vercel. Importar repo lucasbustosmartin-coder/fornitalia. Deploy con preset Other, sin build. App publicada en https://fornitalia.vercel.app/", "Diagnostico"),
    @("27/02/2025", "14:20", "Raíz Vercel con vercel.json", "Crear vercel.json con rewrite: source / → destination /dashboard-flujo-caja.html. Así https://fornitalia.vercel.app/ abre directo el dashboard. Commit y push; Vercel redepliega automático.", "Diagnostico")
)

$r = 12
foreach ($row in $logRows) {
    $log.Cells.Item($r, 1).Value = $row[0]
    $log.Cells.Item($r, 2).Value = $row[1]
    $log.Cells.Item($r, 3).Value = $row[2]
    $log.Cells.Item($r, 4).Value = $row[3]
    $log.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# ---------------------------------------------------------------
# 2. Sheet "Resumen" - append rows 20-21
# ---------------------------------------------------------------
$resumen = $wb.Worksheets.Item("Resumen")

$resumenRows = @(
    @("Repositorio Git (GitHub)", 'Repo: https://github.com/lucasbustosmartin-coder/fornitalia. Rama main. .gitignore excluye node_modules, .venv, .env. Para actualizar: git add . ; git commit -m "mensaje" ; git push origin main.'),
    @("App en producción (Vercel)", "URL pública: https://fornitalia.vercel.app/ (vercel.json reescribe / al dashboard). Cada push a main en GitHub dispara redeploy automático en Vercel. Proyecto: fornitalia, equipo Lucas Bustos, plan Hobby.")
)

$r = 20
foreach ($row in $resumenRows) {
    $resumen.Cells.Item($r, 1).Value = $row[0]
    $resumen.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# ---------------------------------------------------------------
# 3. New sheet "Ref Git y Vercel" (placed after "Resumen", at the end)
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ref = $wb.Worksheets.Add($null, $lastSheet)
$ref.Name = "Ref Git y Vercel"

$refRows = @(
    @("Concepto", "Valor"),
    @("Repositorio GitHub", "https://github.com/lucasbustosmartin-coder/fornitalia"),
    @("URL app en vivo (Vercel)", "https://fornitalia.vercel.app/"),
    @("Rama principal", "main"),
    @("Actualizar y subir cambios", 'git add .  →  git commit -m "descripción"  →  git push origin main'),
    @("Vercel redeploy", "Automático al hacer push a main"),
    @("Archivo configuración raíz", "vercel.json (rewrite / a dashboard-flujo-caja.html)"),
    @("Cuenta GitHub", "lucasbustosmartin-coder"),
    @("Proyecto Vercel", "fornitalia (equipo Lucas Bustos, plan Hobby)")
)

$r = 1
foreach ($row in $refRows) {
    $ref.Cells.Item($r, 1).Value = $row[0]
    $ref.Cells.Item($r, 2).Value = $row[1]
    $r++
}

$ref.Columns.Item(1).ColumnWidth = 28
$ref.Columns.Item(2).ColumnWidth = 70
